$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "331.96"
Set-TextValue "E2" "0.93%"
Set-TextValue "D3" "41.26"
Set-TextValue "E3" "1.71%"
Set-TextValue "D4" "5.686"
Set-TextValue "E4" "-5.84%"
Set-TextValue "D5" "0.08077"
Set-TextValue "E5" "-0.86%"
Set-TextValue "D6" "2.038"
Set-TextValue "E6" "3.51%"
Set-TextValue "D7" "8.745"
Set-TextValue "E7" "-0.28%"
Set-TextValue "D8" "4.542"
Set-TextValue "E8" "-1.43%"
Set-TextValue "E9" "1.97%"
Set-TextValue "D10" "0.9232"
Set-TextValue "E10" "-2.73%"
Set-TextValue "D11" "0.1259"
Set-TextValue "E11" "-5.68%"
Set-TextValue "D12" "0.1949"
Set-TextValue "E12" "-2.36%"
Set-TextValue "D13" "8.729"
Set-TextValue "E13" "-9.22%"
Set-TextValue "D14" "0.09344"
Set-TextValue "E14" "-0.33%"
Set-TextValue "D15" "0.03755"
Set-TextValue "E15" "7.61%"
Set-TextValue "D16" "0.1053"
Set-TextValue "D17" "0.001302"
Set-TextValue "E17" "-0.99%"
Set-TextValue "D18" "0.006267"
Set-TextValue "E18" "0.62%"
Set-TextValue "E19" "0.38%"
Set-TextValue "E20" "-1.80%"
Set-TextValue "D21" "0.1418"
Set-TextValue "E21" "0.26%"
Set-TextValue "D22" "0.2658"
Set-TextValue "E22" "8.60%"
Set-TextValue "D23" "0.04427"
Set-TextValue "E23" "-0.26%"
Set-TextValue "D24" "0.001266"
Set-TextValue "E24" "0.27%"
Set-TextValue "D25" "0.004319"
Set-TextValue "E25" "-1.48%"
Set-TextValue "D26" "0.0001244"
Set-TextValue "E26" "13.82%"
Set-TextValue "D39" "0.02854"
Set-TextValue "E39" "14.91%"
Set-TextValue "D40" "0.05486"
Set-TextValue "E40" "3.62%"
Set-TextValue "D41" "0.007778"
Set-TextValue "E41" "4.08%"
Set-TextValue "D42" "0.009986"
Set-TextValue "E42" "9.80%"
Set-TextValue "E43" "-1.29%"
Set-TextValue "D44" "0.002245"
Set-TextValue "E44" "9.16%"
Set-TextValue "D45" "0.01183"
Set-TextValue "E45" "12.05%"
Set-TextValue "D46" "0.00006776"
Set-TextValue "E46" "-0.74%"
Set-TextValue "D47" "0.00000000753"
Set-TextValue "E47" "0.03%"
Set-TextValue "D48" "0.002287"
Set-TextValue "E48" "26.64%"
Set-TextValue "D49" "0.003017"
Set-TextValue "E49" "-13.77%"
Set-TextValue "D50" "0.00002107"
Set-TextValue "E50" "0.03%"
Set-TextValue "D51" "0.0002007"
Set-TextValue "E51" "0.03%"
